# Insert a new weekly price record at row 430 of the "Ají" (chili pepper)
# sheet for Feria Lagunitas de Puerto Montt, pushing the existing rows
# 430:473 down to 431:474 (dimension grows from A1:R473 to A1:R474).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 430 and below down by one row.
$ws.Rows("430:430").Insert()

# Populate the newly-opened row 430 with the new weekly observation.
$ws.Range("A430").Value = 4
$ws.Range("B430").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C430").Value = "Los Lagos"
$ws.Range("D430").Value = 45212
$ws.Range("E430").Value = 10
$ws.Range("F430").Value = 100112021
$ws.Range("G430").Value = "Ají"
$ws.Range("H430").Value = "Inferno"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 180
$ws.Range("K430").Value = 40000
$ws.Range("L430").Value = 40000
$ws.Range("M430").Value = 40000
$ws.Range("N430").Value = "$/caja 10 kilos"
$ws.Range("O430").Value = "Región de Arica y Parinacota"
$ws.Range("P430").Value = 4000
$ws.Range("Q430").Value = 10
$ws.Range("R430").Value = "Hortaliza"
